$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = 1489253.948956447
$ws.Range("B7").Value = 8112558.798638826
$ws.Range("B10").Value = 5790007.414469049

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("C2").Value = 778312.5217811861
$ws.Range("D2").Value = 778312.5217811862
$ws.Range("E2").Value = 408914.7838453704
$ws.Range("F2").Value = 408914.7838453703
$ws.Range("H2").Value = 438452.9154153269
$ws.Range("I2").Value = 438452.9154153267
$ws.Range("J2").Value = 408914.7838453708
$ws.Range("K2").Value = 408914.7838453705
$ws.Range("L2").Value = 458282.5866920232
$ws.Range("M2").Value = 458282.5866920232
$ws.Range("N2").Value = 458282.5866920232
$ws.Range("O2").Value = 288974.2797907819
$ws.Range("P2").Value = 288974.2797907824
$ws.Range("F4").Value = 246440.2967349876
$ws.Range("H4").Value = 264157.660471655
$ws.Range("I4").Value = 264157.660471655
$ws.Range("K4").Value = 246440.2967349877
$ws.Range("B6").Value = 277024.8260707739
$ws.Range("C6").Value = 277024.8260707738
$ws.Range("D6").Value = 277024.8260707739
$ws.Range("E6").Value = 73676.56641980985
$ws.Range("F6").Value = 145193.8727009419
$ws.Range("G6").Value = 145798.0455821762
$ws.Range("H6").Value = 156651.4194286935
$ws.Range("I6").Value = 156651.4194286933
$ws.Range("J6").Value = 84529.94026632741
$ws.Range("K6").Value = 145193.872700942
$ws.Range("L6").Value = 135206.856271338
$ws.Range("M6").Value = 164331.2792259202
$ws.Range("N6").Value = 164331.2792259202
$ws.Range("O6").Value = 98518.73863735203
$ws.Range("P6").Value = 98518.73863735252

$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("G2").Value = 13.56671730814645

$ws = $wb.Worksheets.Item("Retired Capacities")
$ws.Range("L2").Value = 13.56671730814645

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("J11").Value = 89.39663285141508
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 89.39663285141508
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 65.78242794726917
$ws.Range("Q12").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 65.78242794726972
$ws.Range("Q15").Value = 89.39663285141508
$ws.Range("M17").Value = 102.9633501595615
$ws.Range("K18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("R18").Value = 100.1578341526431
$ws.Range("O19").Value = 88.25797151378666
$ws.Range("K20").Value = 102.9633501595615
$ws.Range("M20").Value = 31.10638119690788
$ws.Range("O20").Value = 102.9633501595615
$ws.Range("P20").Value = 102.9633501595615
$ws.Range("I21").Value = 89.39663285141508
$ws.Range("J21").Value = 0.7465913262578567
$ws.Range("N21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("R21").Value = 100.1578341526431
$ws.Range("J22").Value = 93.35918011667277
$ws.Range("K22").Value = 22.26949182588285
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 2.721440735106512
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 102.9633501595615
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 102.9633501595615
$ws.Range("N25").Value = 60.99283403365629
$ws.Range("J26").Value = 76.81821214870462
$ws.Range("Q26").Value = 89.39663285141508
$ws.Range("K27").Value = 0
$ws.Range("Q27").Value = 89.39663285141508
$ws.Range("J30").Value = 89.39663285141508
$ws.Range("L30").Value = 65.78242794726972
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("L31").Value = 89.39663285141508
$ws.Range("M31").Value = 89.39663285141508
$ws.Range("N31").Value = 89.39663285141508
$ws.Range("O31").Value = 89.39663285141508
$ws.Range("J33").Value = 0.7465913262578567
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("O33").Value = 112.2354442364965
$ws.Range("Q33").Value = 112.2354442364965
$ws.Range("R33").Value = 109.3783679079544
$ws.Range("L34").Value = 112.2354442364965
$ws.Range("O34").Value = 112.2354442364965
$ws.Range("P35").Value = 112.2354442364965
$ws.Range("I36").Value = 89.39663285141508
$ws.Range("K36").Value = 110.6702043477839
$ws.Range("L36").Value = 112.2354442364965
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("R36").Value = 100.1578341526431
$ws.Range("P37").Value = 2.721440735106512
$ws.Range("K39").Value = 107.3790200792832
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 112.2354442364965
$ws.Range("O39").Value = 0
$ws.Range("L40").Value = 112.2354442364965
$ws.Range("M40").Value = 112.2354442364965
$ws.Range("N40").Value = 112.2354442364965
$ws.Range("O40").Value = 112.2354442364965
$ws.Range("J41").Value = 36.40552869322784
$ws.Range("K42").Value = 36.40552869322784
$ws.Range("P42").Value = 0
$ws.Range("K43").Value = 36.4055286932285
$ws.Range("J44").Value = 36.40552869322848
$ws.Range("K45").Value = 36.40552869322784
$ws.Range("L45").Value = 36.40552869322784
$ws.Range("M45").Value = 0
$ws.Range("O45").Value = 36.40552869322784
$ws.Range("Q45").Value = 0

$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("J11").Value = 91.64927167527121
$ws.Range("L12").Value = 138.5543797798742
$ws.Range("M12").Value = 52.73740107060324
$ws.Range("O12").Value = 142.5962444444444
$ws.Range("P12").Value = 68.19197946706107
$ws.Range("Q12").Value = 139.9817740860215
$ws.Range("K15").Value = 137.841438974359
$ws.Range("L15").Value = 72.77195183260446
$ws.Range("Q15").Value = 50.58514123460644
$ws.Range("M17").Value = 127.3828830677112
$ws.Range("K18").Value = 137.841438974359
$ws.Range("N18").Value = 131.3417120833333
$ws.Range("R18").Value = 45.52166981132082
$ws.Range("O19").Value = 50.19856693805613
$ws.Range("K20").Value = 117.126500885419
$ws.Range("M20").Value = 199.2398520303649
$ws.Range("O20").Value = 127.1348612621252
$ws.Range("P20").Value = 128.269645595708
$ws.Range("I21").Value = 10.12574714858493
$ws.Range("J21").Value = 126.0910353404088
$ws.Range("N21").Value = 131.3417120833333
$ws.Range("P21").Value = 133.9744074143302
$ws.Range("R21").Value = 45.52166981132082
$ws.Range("J22").Value = 33.63624132272333
$ws.Range("K22").Value = 106.7437663446525
$ws.Range("L22").Value = 134.8846762812383
$ws.Range("M22").Value = 138.9257839476051
$ws.Range("N22").Value = 127.6855444652332
$ws.Range("O22").Value = 138.4565384518428
$ws.Range("P22").Value = 135.0065633140411
$ws.Range("L24").Value = 138.5543797798742
$ws.Range("M24").Value = 39.17068376245679
$ws.Range("O24").Value = 142.5962444444444
$ws.Range("P24").Value = 31.01105725476872
$ws.Range("N25").Value = 66.6927104315769
$ws.Range("J26").Value = 104.2276923779817
$ws.Range("Q26").Value = 132.9090570230344
$ws.Range("K27").Value = 137.841438974359
$ws.Range("Q27").Value = 50.58514123460644
$ws.Range("J30").Value = 37.44099381525162
$ws.Range("L30").Value = 72.77195183260446
$ws.Range("O30").Value = 142.5962444444444
$ws.Range("P30").Value = 133.9744074143302
$ws.Range("Q30").Value = 139.9817740860215
$ws.Range("L31").Value = 45.48804342982321
$ws.Range("M31").Value = 49.52915109618998
$ws.Range("N31").Value = 38.28891161381812
$ws.Range("O31").Value = 49.05990560042771
$ws.Range("J33").Value = 126.0910353404088
$ws.Range("L33").Value = 138.5543797798742
$ws.Range("M33").Value = 142.1340339220183
$ws.Range("O33").Value = 30.36080020794797
$ws.Range("Q33").Value = 27.74632984952505
$ws.Range("R33").Value = 36.30113605600951
$ws.Range("L34").Value = 22.64923204474182
$ws.Range("O34").Value = 26.22109421534633
$ws.Range("P35").Value = 118.9975515187731
$ws.Range("I36").Value = 10.12574714858493
$ws.Range("K36").Value = 27.17123462657506
$ws.Range("L36").Value = 26.31893554337772
$ws.Range("N36").Value = 131.3417120833333
$ws.Range("O36").Value = 142.5962444444444
$ws.Range("P36").Value = 133.9744074143302
$ws.Range("R36").Value = 45.52166981132082
$ws.Range("P37").Value = 135.0065633140411
$ws.Range("K39").Value = 30.46241889507581
$ws.Range("L39").Value = 138.5543797798742
$ws.Range("M39").Value = 142.1340339220183
$ws.Range("N39").Value = 19.10626784683684
$ws.Range("O39").Value = 142.5962444444444
$ws.Range("L40").Value = 22.64923204474182
$ws.Range("M40").Value = 26.69033971110859
$ws.Range("N40").Value = 15.45010022873673
$ws.Range("O40").Value = 26.22109421534633
$ws.Range("J41").Value = 144.6403758334584
$ws.Range("K42").Value = 101.4359102811312
$ws.Range("P42").Value = 133.9744074143302
$ws.Range("K43").Value = 92.60772947730688
$ws.Range("J44").Value = 144.6403758334578
$ws.Range("K45").Value = 101.4359102811312
$ws.Range("L45").Value = 102.1488510866463
$ws.Range("M45").Value = 142.1340339220183
$ws.Range("O45").Value = 106.1907157512166
$ws.Range("Q45").Value = 139.9817740860215

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B7").Value = 468642.4575694173
$ws.Range("B8").Value = 468642.4575694173
$ws.Range("B9").Value = 468642.4575694173
$ws.Range("B12").Value = 482398.1654581949
$ws.Range("B13").Value = 482398.1654581949

Write-Output "Applied 211 cell updates"